$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated FilesTab query (B4): dropped the `File Type` and `Breed` output
# columns from the RETURN clause (ICDC Breed script correction).
$filesTabQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Basset Hound']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $filesTabQuery

# Row 4 shrank now that the query has two fewer output columns listed.
$ws.Rows.Item(4).RowHeight = 217.5

# Selection / scroll moved from row 3 down onto the FilesTab query cell.
$ws.Range("B4").Select()
